# "updated legacy GSC export data"
# The old export had a leading placeholder row for 2025-10-10 (with blank
# Not-indexed/Indexed values) that the refreshed export no longer contains.
# Remove that row from the "Chart" sheet; every later row (and the
# sharedStrings table) shifts up/renumbers accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$ws.Rows.Item(2).Delete() | Out-Null
